# AuthorList.xlsx edit: add Weizmann footnote for R. Felkai, add S.R. Soleti at DIPC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) R. Felkai (row 31): add "Now at Weizmann Institute of Science, Israel." footnote ---
# (Do this first so the new shared string lands at the same index the
# canonical file uses: 246.)
$ws.Range("C31").Value = " Now at Weizmann Institute of Science, Israel."

# --- 2) Add S.R. Soleti (DIPC) as a new row, right before row 94 (Sorel) ---
# Grab the existing DIPC institution/address text (row 98: Torrent) so the
# new row reuses the same shared strings instead of creating near-duplicates.
$inst = $ws.Range("D98").Value()
$addr = $ws.Range("E98").Value()

$ws.Rows(94).Insert()

$ws.Range("A94").Value = "Soleti"
$ws.Range("B94").Value = "S.R."
$ws.Range("D94").Value = $inst
$ws.Range("E94").Value = $addr

# Match the "black" font override used by the other DIPC address cell
# (row 99 after the insert, same row that D98/E98 came from) rather than
# the plain default formatting.
$ws.Range("D99").Copy()
$ws.Range("D94").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3) Reflect the cursor/scroll state recorded in the saved file ---
$ws.Activate()
$ws.Range("E94").Select()
